# Weekly price update: insert a new "Fruta / hortaliza" record at the top of
# the Mango / Macroferia Regional de Talca price history (row 157), pushing
# all existing rows (157-180) down by one (to 158-181).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 157:180 down to 158:181, leaving a blank row 157 to fill in.
$ws.Rows("157:157").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A157").Value = 5
$ws.Range("B157").Value = "Macroferia Regional de Talca"
$ws.Range("C157").Value = "Maule"
$ws.Range("D157").Value = 45131
$ws.Range("E157").Value = 7
$ws.Range("F157").Value = "Fruta"
$ws.Range("G157").Value = 100108
$ws.Range("H157").Value = "Tropicales y subtropicales"
$ws.Range("I157").Value = 100108002
$ws.Range("J157").Value = "Mango"
$ws.Range("K157").Value = "Sin especificar"
$ws.Range("L157").Value = "Primera"
$ws.Range("M157").Value = 248
$ws.Range("N157").Value = 8000
$ws.Range("O157").Value = 8000
$ws.Range("P157").Value = 8000
$ws.Range("Q157").Value = "$/bandeja 4 kilos"
$ws.Range("R157").Value = "Brasil"
$ws.Range("S157").Value = 2000
$ws.Range("T157").Value = 4
